# Auto-generated edit script applying numeric market-data updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve-profit worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 132.5
$ws.Range("I2").Value = 132.5
$ws.Range("K2").Value = 132.5
$ws.Range("M2").Value = -19.5
$ws.Range("H76").Value = 100006600
$ws.Range("J76").Value = 125006500
$ws.Range("L76").Value = 125006500
$ws.Range("N76").Value = -125007130
$ws.Range("H79").Value = 100006600
$ws.Range("J79").Value = 125006500
$ws.Range("L79").Value = 125006500
$ws.Range("N79").Value = -125008684
$ws.Range("H111").Value = 112722.555
$ws.Range("J111").Value = 300
$ws.Range("L111").Value = 900
$ws.Range("N111").Value = -7034
$ws.Range("H129").Value = 2296.1667
$ws.Range("I129").Value = 731.375
$ws.Range("J129").Value = 3078.5625
$ws.Range("K129").Value = 2194.125
$ws.Range("L129").Value = 9235.6875
$ws.Range("M129").Value = 2805.875
$ws.Range("N129").Value = -19235.6875
$ws.Range("H132").Value = 1384.6271
$ws.Range("I132").Value = 1343.8246
$ws.Range("K132").Value = 4031.4738
$ws.Range("M132").Value = -1501.4738
$ws.Range("H134").Value = 77322.60000000001
$ws.Range("J134").Value = 77322.60000000001
$ws.Range("L134").Value = 77322.60000000001
$ws.Range("N134").Value = -87462.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4071.2043
$ws.Range("I32").Value = 3054.3765
$ws.Range("K32").Value = 3054.3765
$ws.Range("M32").Value = -2767.3765
$ws.Range("H74").Value = 2502.348
$ws.Range("I74").Value = 2478.375
$ws.Range("K74").Value = 2478.375
$ws.Range("M74").Value = -1604.375
$ws.Range("H77").Value = 2502.348
$ws.Range("I77").Value = 2478.375
$ws.Range("K77").Value = 12391.875
$ws.Range("M77").Value = -8023.875
$ws.Range("H110").Value = 168187.5
$ws.Range("I110").Value = 209834.5
$ws.Range("K110").Value = 209834.5
$ws.Range("M110").Value = -207789.5
$ws.Range("H122").Value = 5528.143
$ws.Range("I122").Value = 3982.5
$ws.Range("K122").Value = 11947.5
$ws.Range("M122").Value = -9497.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 305020.97
$ws.Range("I107").Value = 1851
$ws.Range("K107").Value = 1851
$ws.Range("M107").Value = 69

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 52262.4
$ws.Range("I31").Value = 1376.6
$ws.Range("J31").Value = 204919.8
$ws.Range("K31").Value = 1376.6
$ws.Range("L31").Value = 204919.8
$ws.Range("M31").Value = -1081.6
$ws.Range("N31").Value = -205509.8
$ws.Range("H34").Value = 52262.4
$ws.Range("I34").Value = 1376.6
$ws.Range("J34").Value = 204919.8
$ws.Range("K34").Value = 1376.6
$ws.Range("L34").Value = 204919.8
$ws.Range("M34").Value = -1174.6
$ws.Range("N34").Value = -205323.8
$ws.Range("H58").Value = 2423
$ws.Range("I58").Value = 1995.5
$ws.Range("K58").Value = 1995.5
$ws.Range("M58").Value = -1792.5
$ws.Range("H59").Value = 41639.125
$ws.Range("J59").Value = 41639.125
$ws.Range("L59").Value = 41639.125
$ws.Range("N59").Value = -43929.125
$ws.Range("H86").Value = 7815.4287
$ws.Range("I86").Value = 6000
$ws.Range("K86").Value = 6000
$ws.Range("M86").Value = -4877
$ws.Range("H89").Value = 7815.4287
$ws.Range("I89").Value = 6000
$ws.Range("K89").Value = 30000
$ws.Range("M89").Value = -24384
$ws.Range("H122").Value = 2680.6667
$ws.Range("I122").Value = 2515.9375
$ws.Range("K122").Value = 7547.8125
$ws.Range("M122").Value = -5097.8125
$ws.Range("H136").Value = 2423
$ws.Range("I136").Value = 1995.5
$ws.Range("K136").Value = 5986.5
$ws.Range("M136").Value = -3436.5
$ws.Range("H141").Value = 221417.05
$ws.Range("J141").Value = 220924.58
$ws.Range("L141").Value = 220924.58
$ws.Range("N141").Value = -231284.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 19609238
$ws.Range("J129").Value = 2504.4285
$ws.Range("L129").Value = 7513.2855
$ws.Range("N129").Value = -17513.2855
$ws.Range("H132").Value = 922060.7
$ws.Range("I132").Value = 205167.4
$ws.Range("K132").Value = 1846506.6
$ws.Range("M132").Value = -1843976.6
$ws.Range("H139").Value = 5738.95
$ws.Range("I139").Value = 2726.762
$ws.Range("K139").Value = 8180.286
$ws.Range("M139").Value = -3040.286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 1000000000
$ws.Range("J21").Value = 1000000000
$ws.Range("L21").Value = 1000000000
$ws.Range("N21").Value = -1000000346
$ws.Range("H24").Value = 1373862.9
$ws.Range("I24").Value = 27666.666
$ws.Range("J24").Value = 1741007.4
$ws.Range("K24").Value = 27666.666
$ws.Range("L24").Value = 1741007.4
$ws.Range("M24").Value = -27493.666
$ws.Range("N24").Value = -1741353.4
$ws.Range("H30").Value = 1000000000
$ws.Range("J30").Value = 1000000000
$ws.Range("L30").Value = 1000000000
$ws.Range("N30").Value = -1000000210
$ws.Range("H36").Value = 1950
$ws.Range("J36").Value = 2200
$ws.Range("L36").Value = 2200
$ws.Range("N36").Value = -3170
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 110140
$ws.Range("J86").Value = 110140
$ws.Range("L86").Value = 110140
$ws.Range("N86").Value = -112512
$ws.Range("H89").Value = 110140
$ws.Range("J89").Value = 110140
$ws.Range("L89").Value = 330420
$ws.Range("N89").Value = -342276
$ws.Range("H97").Value = 579.3333
$ws.Range("I97").Value = 631
$ws.Range("K97").Value = 631
$ws.Range("M97").Value = -135
$ws.Range("H102").Value = 3143.4062
$ws.Range("I102").Value = 1569.1305
$ws.Range("K102").Value = 1569.1305
$ws.Range("M102").Value = 52.86950000000002
$ws.Range("H113").Value = 597548.1
$ws.Range("I113").Value = 1114202.4
$ws.Range("K113").Value = 1114202.4
$ws.Range("M113").Value = -1112032.4
$ws.Range("H126").Value = 3252
$ws.Range("I126").Value = 2279.7
$ws.Range("K126").Value = 6839.099999999999
$ws.Range("M126").Value = -4369.099999999999
$ws.Range("H138").Value = 45065
$ws.Range("I138").Value = 20390
$ws.Range("J138").Value = 50000
$ws.Range("K138").Value = 20390
$ws.Range("L138").Value = 50000
$ws.Range("M138").Value = -15250
$ws.Range("N138").Value = -60280
$ws.Range("H140").Value = 132500
$ws.Range("J140").Value = 132500
$ws.Range("L140").Value = 132500
$ws.Range("N140").Value = -142860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4564.5713
$ws.Range("I7").Value = 4318.2666
$ws.Range("K7").Value = 4318.2666
$ws.Range("M7").Value = -4206.2666
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H122").Value = 505112.34
$ws.Range("I122").Value = 836520.5600000001
$ws.Range("K122").Value = 2509561.68
$ws.Range("M122").Value = -2507111.68
$ws.Range("H126").Value = 4564.5713
$ws.Range("I126").Value = 4318.2666
$ws.Range("K126").Value = 12954.7998
$ws.Range("M126").Value = -10484.7998
$ws.Range("H136").Value = 381229.12
$ws.Range("I136").Value = 564899.2
$ws.Range("J136").Value = 13889
$ws.Range("K136").Value = 1694697.6
$ws.Range("L136").Value = 41667
$ws.Range("M136").Value = -1692147.6
$ws.Range("N136").Value = -46767

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 6949.7
$ws.Range("J33").Value = 7499.6665
$ws.Range("L33").Value = 7499.6665
$ws.Range("N33").Value = -7999.6665
$ws.Range("H36").Value = 6949.7
$ws.Range("J36").Value = 7499.6665
$ws.Range("L36").Value = 7499.6665
$ws.Range("N36").Value = -7999.6665
$ws.Range("H54").Value = 29019.25
$ws.Range("J54").Value = 30692.334
$ws.Range("L54").Value = 30692.334
$ws.Range("N54").Value = -31732.334
$ws.Range("H122").Value = 27029598
$ws.Range("I122").Value = 47620060
$ws.Range("K122").Value = 142860180
$ws.Range("M122").Value = -142857730
$ws.Range("H125").Value = 56666.668
$ws.Range("J125").Value = 56666.668
$ws.Range("L125").Value = 56666.668
$ws.Range("N125").Value = -66506.66800000001
$ws.Range("H126").Value = 1576.6666
$ws.Range("I126").Value = 1565.24
$ws.Range("J126").Value = 1719.5
$ws.Range("K126").Value = 4695.72
$ws.Range("L126").Value = 5158.5
$ws.Range("M126").Value = -2225.72
$ws.Range("N126").Value = -10098.5
$ws.Range("H141").Value = 58833
$ws.Range("J141").Value = 58833
$ws.Range("L141").Value = 58833
$ws.Range("N141").Value = -69193
